$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: Cloud 291 ("intranet") ---
# Reposition/resize the cloud and prefix its text with two spaces.
$cloud = $s.Shapes.Item(1)
$cloud.Top = 2917818 / 12700
$cloud.Height = 543310 / 12700
$cloudText = $cloud.TextFrame.TextRange
$cloudText.Text = "  intranet"

# --- Shape 5: Rectangle 42 ("Core and Middle Tier Modules") ---
# Drop "and Middle Tier " from the middle, leaving "Core " + "Modules" as
# two separate runs.
$coreRect = $s.Shapes.Item(5)
$coreText = $coreRect.TextFrame.TextRange
$coreMid = $coreText.Characters(6, 16)
$coreMid.Text = ""

# --- Shape 6: Rectangle 10 ("Voice command" / "input listener") ---
# Rename "Voice command" to "Speech command" (merging the first two runs)
# and merge "input " + "listener" into a single run, keeping the <a:br/>.
$voiceRect = $s.Shapes.Item(6)
$voiceText = $voiceRect.TextFrame.TextRange
$voiceFirst = $voiceText.Characters(1, 13)
$voiceFirst.Text = "Speech command"
$voiceText2 = $voiceRect.TextFrame.TextRange
$voiceSecond = $voiceText2.Characters(16, 14)
$voiceSecond.Text = "input listener"

# --- Shape 8: Rectangle 69 ("Laser/mouse " / "input " / "listener") ---
# Merge the three runs into a single run. Go through a Characters() sub-range
# (rather than TextRange.Text directly) so the merge happens even though the
# concatenated text reads the same before and after.
$laserRect = $s.Shapes.Item(8)
$laserText = $laserRect.TextFrame.TextRange
$laserAll = $laserText.Characters(1, $laserText.Text.Length)
$laserAll.Text = "Laser/mouse input listener"

# --- Shape 9: Rectangle 70 ("Tracking postures " / "input " / "listener") ---
# Merge the three runs into a single run (same approach as above).
$trackRect = $s.Shapes.Item(9)
$trackText = $trackRect.TextFrame.TextRange
$trackAll = $trackText.Characters(1, $trackText.Text.Length)
$trackAll.Text = "Tracking postures input listener"
